$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.447159934046681
$ws.Range("C2").Value = 0.1183654977517392
$ws.Range("D2").Value = 0.05647233328883239
$ws.Range("E2").Value = 0.1211252422109226
$ws.Range("F2").Value = 1.450653168606735
$ws.Range("I2").Value = 1.012805003001837
$ws.Range("K2").Value = 0.5539443074843007
$ws.Range("M2").Value = 0.2907387100251597
$ws.Range("N2").Value = 2.175246244320618
$ws.Range("B3").Value = 0.4095024174870048
$ws.Range("C3").Value = 0.1075380643835331
$ws.Range("D3").Value = 0.05652549211420244
$ws.Range("E3").Value = 0.1114079124986631
$ws.Range("F3").Value = 1.432167997306735
$ws.Range("I3").Value = 1.008086454538415
$ws.Range("K3").Value = 0.5062130379695304
$ws.Range("M3").Value = 0.2664624142609497
$ws.Range("N3").Value = 2.185885919413586
$ws.Range("B4").Value = 0.3865995802258908
$ws.Range("C4").Value = 0.1009493183563563
$ws.Range("D4").Value = 0.05655676561888257
$ws.Range("E4").Value = 0.105515144135552
$ws.Range("F4").Value = 1.421597236561666
$ws.Range("I4").Value = 1.005664435391168
$ws.Range("K4").Value = 0.4771790985441271
$ws.Range("M4").Value = 0.2517154666154795
$ws.Range("N4").Value = 2.193064949504098
$ws.Range("B5").Value = 0.3773215009832143
$ws.Range("C5").Value = 0.09827912791155313
$ws.Range("D5").Value = 0.05656917047832444
$ws.Range("E5").Value = 0.103132072231908
$ws.Range("F5").Value = 1.417485132810867
$ws.Range("I5").Value = 1.004796747489628
$ws.Range("K5").Value = 0.4654160143265642
$ws.Range("M5").Value = 0.2457455522015977
$ws.Range("N5").Value = 2.196152724653381
$ws.Range("B6").Value = 0.3757842039003663
$ws.Range("C6").Value = 0.09783663311780799
$ws.Range("D6").Value = 0.05657120992673725
$ws.Range("E6").Value = 0.1027374615083048
$ws.Range("F6").Value = 1.416814122139456
$ws.Range("I6").Value = 1.004659868909279
$ws.Range("K6").Value = 0.4634668927290591
$ws.Range("M6").Value = 0.244756633405359
$ws.Range("N6").Value = 2.196675240725632
$ws.Range("B7").Value = 0.3864742301406352
$ws.Range("C7").Value = 0.1009132475935388
$ws.Range("D7").Value = 0.05655693428415631
$ws.Range("E7").Value = 0.1054829315197665
$ws.Range("F7").Value = 1.421540987946074
$ws.Range("I7").Value = 1.005652250639848
$ws.Range("K7").Value = 0.4770201804957424
$ws.Range("M7").Value = 0.2516347944289876
$ws.Range("N7").Value = 2.193105935653193
$ws.Range("B8").Value = 0.4341301106048832
$ws.Range("C8").Value = 0.1146197913980131
$ws.Range("D8").Value = 0.05649094849302827
$ws.Range("E8").Value = 0.1177592505251539
$ws.Range("F8").Value = 1.444117522820079
$ws.Range("I8").Value = 1.011079306099063
$ws.Range("K8").Value = 0.5374296969102943
$ws.Range("M8").Value = 0.2823350696117899
$ws.Range("N8").Value = 2.178780536494799
$ws.Range("B9").Value = 0.5293293384286244
$ws.Range("C9").Value = 0.1419771266247949
$ws.Range("D9").Value = 0.05635052502602633
$ws.Range("E9").Value = 0.1424307895993593
$ws.Range("F9").Value = 1.49459365715083
$ws.Range("I9").Value = 1.025502279070807
$ws.Range("K9").Value = 0.6580778015864439
$ws.Range("M9").Value = 0.3438170289361366
$ws.Range("N9").Value = 2.155826967318575
$ws.Range("B10").Value = 0.6003559906597786
$ws.Range("C10").Value = 0.1623816424976781
$ws.Range("D10").Value = 0.05624038460955028
$ws.Range("E10").Value = 0.1609409972826086
$ws.Range("F10").Value = 1.535495407624055
$ws.Range("I10").Value = 1.038420108519325
$ws.Range("K10").Value = 0.7480834786603054
$ws.Range("M10").Value = 0.389798283343481
$ws.Range("N10").Value = 2.142111173269996
$ws.Range("B11").Value = 0.6329078513584818
$ws.Range("C11").Value = 0.1717334078408044
$ws.Range("D11").Value = 0.05618871546058024
$ws.Range("E11").Value = 0.1694495231298774
$ws.Range("F11").Value = 1.554939478887007
$ws.Range("I11").Value = 1.04480457637279
$ws.Range("K11").Value = 0.7893337248197554
$ws.Range("M11").Value = 0.4108993008220665
$ws.Range("N11").Value = 2.136558583585753
$ws.Range("B12").Value = 0.6452693322060554
$ws.Range("C12").Value = 0.1752848905521773
$ws.Range("D12").Value = 0.0561689210750913
$ws.Range("E12").Value = 0.1726844711147137
$ws.Range("F12").Value = 1.562423431925922
$ws.Range("I12").Value = 1.047295532701796
$ws.Range("K12").Value = 0.8049985816886362
$ws.Range("M12").Value = 0.4189166185275184
$ws.Range("N12").Value = 2.134555017945246
$ws.Range("B13").Value = 0.6426055172895815
$ws.Range("C13").Value = 0.1745195602853471
$ws.Range("D13").Value = 0.05617319435787138
$ws.Range("E13").Value = 0.171987186581454
$ws.Range("F13").Value = 1.560806244203661
$ws.Range("I13").Value = 1.046755797007052
$ws.Range("K13").Value = 0.8016229008925393
$ws.Range("M13").Value = 0.4171887468803703
$ws.Range("N13").Value = 2.134982109923968
$ws.Range("B14").Value = 0.6339241399982996
$ws.Range("C14").Value = 0.1720253856929617
$ws.Range("D14").Value = 0.05618709155979218
$ws.Range("E14").Value = 0.1697154030222663
$ws.Range("F14").Value = 1.555552762486087
$ws.Range("I14").Value = 1.045008038511881
$ws.Range("K14").Value = 0.7906215930505596
$ws.Range("M14").Value = 0.4115583507014762
$ws.Range("N14").Value = 2.136391761487332
$ws.Range("B15").Value = 0.6286110796556557
$ws.Range("C15").Value = 0.1704989611276062
$ws.Range("D15").Value = 0.05619557415973553
$ws.Range("E15").Value = 0.1683255653134026
$ws.Range("F15").Value = 1.552350612142249
$ws.Range("I15").Value = 1.043947038422985
$ws.Range("K15").Value = 0.7838887489484705
$ws.Range("M15").Value = 0.408113074421486
$ws.Range("N15").Value = 2.13726812641454
$ws.Range("B16").Value = 0.5982335153716178
$ws.Range("C16").Value = 0.1617718975122955
$ws.Range("D16").Value = 0.0562437295211069
$ws.Range("E16").Value = 0.1603867425453345
$ws.Range("F16").Value = 1.534241586109516
$ws.Range("I16").Value = 1.038013112491363
$ws.Range("K16").Value = 0.7453938611702426
$ws.Range("M16").Value = 0.3884230149123979
$ws.Range("N16").Value = 2.142487895312769
$ws.Range("B17").Value = 0.5796597633303975
$ws.Range("C17").Value = 0.1564360763944137
$ws.Range("D17").Value = 0.05627286787396812
$ws.Range("E17").Value = 0.1555393047417652
$ws.Range("F17").Value = 1.523347149060896
$ws.Range("I17").Value = 1.03450313970545
$ws.Range("K17").Value = 0.7218570996028006
$ws.Range("M17").Value = 0.3763911707498835
$ws.Range("N17").Value = 2.145866184696715
$ws.Range("B18").Value = 0.568999357085346
$ws.Range("C18").Value = 0.153373599422082
$ws.Range("D18").Value = 0.05628948028958725
$ws.Range("E18").Value = 0.1527594659934195
$ws.Range("F18").Value = 1.517159774224538
$ws.Range("I18").Value = 1.032532105045796
$ws.Range("K18").Value = 0.7083481837056809
$ws.Range("M18").Value = 0.36948806054356
$ws.Range("N18").Value = 2.147873912245288
$ws.Range("B19").Value = 0.565393827520694
$ws.Range("C19").Value = 0.1523378169046055
$ws.Range("D19").Value = 0.05629507979064208
$ws.Range("E19").Value = 0.1518196738464326
$ws.Range("F19").Value = 1.515078358710284
$ws.Range("I19").Value = 1.031872950681176
$ws.Range("K19").Value = 0.7037792369341105
$ws.Range("M19").Value = 0.3671537439954662
$ws.Range("N19").Value = 2.148564784034392
$ws.Range("B20").Value = 0.5816346185085024
$ws.Range("C20").Value = 0.1570034050243407
$ws.Range("D20").Value = 0.05626978130382199
$ws.Range("E20").Value = 0.1560544638830592
$ws.Range("F20").Value = 1.524498719941818
$ws.Range("I20").Value = 1.034871832562715
$ws.Range("K20").Value = 0.7243596448546725
$ws.Range("M20").Value = 0.3776701891325089
$ws.Range("N20").Value = 2.145499869281167
$ws.Range("B21").Value = 0.6364731259482994
$ws.Range("C21").Value = 0.1727577079044806
$ws.Range("D21").Value = 0.05618301583697871
$ws.Range("E21").Value = 0.1703823271179061
$ws.Range("F21").Value = 1.557092552459963
$ws.Range("I21").Value = 1.045519406400658
$ws.Range("K21").Value = 0.793851739961525
$ws.Range("M21").Value = 0.4132114036494912
$ws.Range("N21").Value = 2.135975021091966
$ws.Range("B22").Value = 0.6725161247629785
$ws.Range("C22").Value = 0.1831134204099101
$ws.Range("D22").Value = 0.05612497694883167
$ws.Range("E22").Value = 0.1798220506826311
$ws.Range("F22").Value = 1.579099447477773
$ws.Range("I22").Value = 1.052905553272126
$ws.Range("K22").Value = 0.8395270819161738
$ws.Range("M22").Value = 0.4365961474293414
$ws.Range("N22").Value = 2.130327645275216
$ws.Range("B23").Value = 0.6532607107407955
$ws.Range("C23").Value = 0.1775808976678377
$ws.Range("D23").Value = 0.05615607635110997
$ws.Range("E23").Value = 0.1747768788225343
$ws.Range("F23").Value = 1.567289300948943
$ws.Range("I23").Value = 1.048924247357604
$ws.Range("K23").Value = 0.8151255790032224
$ws.Range("M23").Value = 0.4241008248648797
$ws.Range("N23").Value = 2.133288793090955
$ws.Range("B24").Value = 0.5807417307613321
$ws.Range("C24").Value = 0.1567468996993853
$ws.Range("D24").Value = 0.0562711771759119
$ws.Range("E24").Value = 0.1558215386016499
$ws.Range("F24").Value = 1.523977858088699
$ws.Range("I24").Value = 1.034705000452433
$ws.Range("K24").Value = 0.7232281735190895
$ws.Range("M24").Value = 0.3770919008587299
$ws.Range("N24").Value = 2.145665276436716
$ws.Range("B25").Value = 0.5033865067819079
$ws.Range("C25").Value = 0.134523623584613
$ws.Range("D25").Value = 0.05638972030311606
$ws.Range("E25").Value = 0.1356904722248728
$ws.Range("F25").Value = 1.480270782167182
$ws.Range("I25").Value = 1.021194050073603
$ws.Range("K25").Value = 0.6252019211409561
$ws.Range("M25").Value = 0.3270444468714686
$ws.Range("N25").Value = 2.161484806251352
